$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Locate the anchor paragraph: "...Else they go employee group." is followed
# by two blank paragraphs; the new "Procedures"/"Functions" content is
# inserted right after the second blank paragraph (one blank paragraph is
# left before the final sectPr, matching the target document).
# ---------------------------------------------------------------------------
$anchorIndex = -1
$total = $d.Paragraphs.Count
for ($i = 1; $i -le $total; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*Else they go employee group.*") {
        $anchorIndex = $i
    }
}

$cur = $anchorIndex + 2

# Grab an existing bulleted ("ListParagraph" / numId 1) paragraph so the new
# bullet items can continue the very same numbered/bulleted list instead of
# minting a new list definition.
$bulletTemplate = $d.Paragraphs(7).Range.ListFormat.ListTemplate

# ---------------------------------------------------------------------------
# "Procedures" section
# ---------------------------------------------------------------------------
$d.Paragraphs($cur).Range.InsertParagraphAfter()
$cur = $cur + 1
$d.Paragraphs($cur).Range.Text = "Procedures"
$d.Paragraphs($cur).Style = "Heading 2"

$d.Paragraphs($cur).Range.InsertParagraphAfter()
$cur = $cur + 1
$d.Paragraphs($cur).Range.Text = "Made procedures:"
$d.Paragraphs($cur).Style = "Normal"

$d.Paragraphs($cur).Range.InsertParagraphAfter()
$cur = $cur + 1
$d.Paragraphs($cur).Range.Text = "salaryBase()"
$d.Paragraphs($cur).Style = "List Paragraph"
$d.Paragraphs($cur).Range.ListFormat.ApplyListTemplateWithLevel($bulletTemplate, $true, 1, $false, $false)

$d.Paragraphs($cur).Range.InsertParagraphAfter()
$cur = $cur + 1
$d.Paragraphs($cur).Range.Text = "temporaryIncrease()"
$d.Paragraphs($cur).Style = "List Paragraph"
$d.Paragraphs($cur).Range.ListFormat.ApplyListTemplateWithLevel($bulletTemplate, $true, 1, $false, $false)

$d.Paragraphs($cur).Range.InsertParagraphAfter()
$cur = $cur + 1
$d.Paragraphs($cur).Range.Text = "percentSalaryIncrease(percentValue numeric, maximumValue numeric)"
$d.Paragraphs($cur).Style = "List Paragraph"
$d.Paragraphs($cur).Range.ListFormat.ApplyListTemplateWithLevel($bulletTemplate, $true, 1, $false, $false)

$d.Paragraphs($cur).Range.InsertParagraphAfter()
$cur = $cur + 1
$d.Paragraphs($cur).Range.Text = "correctSalary()"
$d.Paragraphs($cur).Style = "List Paragraph"
$d.Paragraphs($cur).Range.ListFormat.ApplyListTemplateWithLevel($bulletTemplate, $true, 1, $false, $false)

$d.Paragraphs($cur).Range.InsertParagraphAfter()
$cur = $cur + 1
$d.Paragraphs($cur).Style = "Normal"

$d.Paragraphs($cur).Range.InsertParagraphAfter()
$cur = $cur + 1
$d.Paragraphs($cur).Range.Text = "SalaryBase sets all employees salary to the salary given by their job title."
$d.Paragraphs($cur).Style = "Normal"

$d.Paragraphs($cur).Range.InsertParagraphAfter()
$cur = $cur + 1
$d.Paragraphs($cur).Range.Text = "TemporaryIncrease gives all employees with temporary contract 3 months more contract time"
$d.Paragraphs($cur).Style = "Normal"

$d.Paragraphs($cur).Range.InsertParagraphAfter()
$cur = $cur + 1
$d.Paragraphs($cur).Range.Text = "percentSalaryIncrease takes in any numeric values. percentValue proceeds with integers being per cents like 20 = 20% and so on. It also takes maximum value and if the original value was higher than the given value then the value wasn" + [char]8217 + "t increased. percentSalaryIncrease increased current salary value by given per cent value."
$d.Paragraphs($cur).Style = "Normal"

$d.Paragraphs($cur).Range.InsertParagraphAfter()
$cur = $cur + 1
$d.Paragraphs($cur).Range.Text = "correctSalary first call salaryBase procedure to give them their salary a base value and then gives them additional salary for each benefit salary marked in the skills they have."
$d.Paragraphs($cur).Style = "Normal"

$d.Paragraphs($cur).Range.InsertParagraphAfter()
$cur = $cur + 1
$d.Paragraphs($cur).Style = "Normal"

# ---------------------------------------------------------------------------
# "Functions" section
# ---------------------------------------------------------------------------
$d.Paragraphs($cur).Range.InsertParagraphAfter()
$cur = $cur + 1
$d.Paragraphs($cur).Range.Text = "Functions"
$d.Paragraphs($cur).Style = "Heading 2"

$d.Paragraphs($cur).Range.InsertParagraphAfter()
$cur = $cur + 1
$d.Paragraphs($cur).Style = "Normal"

$d.Paragraphs($cur).Range.InsertParagraphAfter()
$cur = $cur + 1
$d.Paragraphs($cur).Range.Text = "Made function:"
$d.Paragraphs($cur).Style = "Normal"

$d.Paragraphs($cur).Range.InsertParagraphAfter()
$cur = $cur + 1
$d.Paragraphs($cur).Range.Text = "getProjects(givenDate date) "
$d.Paragraphs($cur).Style = "List Paragraph"
$d.Paragraphs($cur).Range.ListFormat.ApplyListTemplateWithLevel($bulletTemplate, $true, 1, $false, $false)

$d.Paragraphs($cur).Range.InsertParagraphAfter()
$cur = $cur + 1
$d.Paragraphs($cur).Style = "Normal"

$d.Paragraphs($cur).Range.InsertParagraphAfter()
$cur = $cur + 1
$d.Paragraphs($cur).Range.Text = "getProjects takes a date and returns in a table all projects which end date is later than given date. In the table are information about the project and the customer information."
$d.Paragraphs($cur).Style = "Normal"

Write-Output "done; final cur=$cur ; paragraph count=$($d.Paragraphs.Count)"
